$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.397.90"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "3.073.84"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.00%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "522.49"
$ws.Range("E5").Value = "  +1.38%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "140.16"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.073.18"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.80%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "7.14"
$ws.Range("E10").Value = "  -2.59%  "
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("E12").Value = "  +2.52%  "
$ws.Range("D13").Value = "3.603.08"
$ws.Range("E13").Value = "  +0.03%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.133"
$ws.Range("E14").Value = "  +2.00%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "25.32"
$ws.Range("E15").Value = "  -5.20%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000163"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").Value = "57.480.47"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "3.075.84"
$ws.Range("E18").Value = "  -0.08%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.06"
$ws.Range("E19").Value = "  -1.66%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.62"
$ws.Range("E20").Value = "  -1.60%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "7.91"
$ws.Range("E21").Value = "  -2.23%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "338.43"
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("E23").Value = "  +0.21%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.508"
$ws.Range("E24").Value = "  +1.20%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "66.80"
$ws.Range("E25").Value = "  +2.62%  "
$ws.Range("E26").Value = "  -2.07%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "0.0₃0904"
$ws.Range("E28").Value = "  +0.21%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.00"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  -1.56%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.19"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("E32").Value = "  +3.01%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "20.84"
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("E34").Value = "  -2.31%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "157.80"
$ws.Range("E35").Value = "  +1.83%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "4.57"
$ws.Range("E36").Value = "  +0.66%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "6.10"
$ws.Range("E37").Value = "  +1.75%  "
$ws.Range("E38").Value = "  -5.33%  "
$ws.Range("E39").Value = "  -2.11%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0660"
$ws.Range("E40").Value = "  -2.04%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.56"
$ws.Range("E41").Value = "  +12.59%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.97"
$ws.Range("E42").Value = "  +2.01%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.680"
$ws.Range("E43").Value = "  +3.76%  "
$ws.Range("D44").Value = "3.112.10"
$ws.Range("E44").Value = "  -0.13%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "36.71"
$ws.Range("E45").Value = "  +0.36%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.00"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "2.269.97"
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("E48").Value = "  +2.40%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.986"
$ws.Range("E49").Value = "  +4.70%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "6.05"
$ws.Range("E50").Value = "  +1.75%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "20.32"
$ws.Range("E51").Value = "  -0.77%  "
